$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Sheet "Sheet1": update OrderDate (A2) and OverageID (L2).
# Both cells are stored as plain text (shared-string) in the workbook, but
# naively assigning a date-looking / number-looking string via .Value
# causes Excel's auto-detection to coerce it into a date/number and fork a
# brand-new cell style. To keep the values as text (matching the original
# file) and to keep the existing cell style (fill/border) untouched, stage
# the new text in a scratch cell that has already been forced to Text
# format, copy it, and paste-special (values only) into the target cell.
# The scratch cell is then fully deleted (with a left-shift) so it leaves
# no trace in the sheet's dimension or style table.
# -----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# OrderDate: 10-07-2021 -> 10-15-2021
$ws1.Range("N1").NumberFormat = "@"
$ws1.Range("N1").Value = "10-15-2021"
$ws1.Range("N1").Copy()
$ws1.Range("A2").PasteSpecial(-4163)   # xlPasteValues
$ws1.Range("N1").Delete(-4159)         # xlShiftToLeft

# OverageID: 57762255 -> 57794961
$ws1.Range("N1").NumberFormat = "@"
$ws1.Range("N1").Value = "57794961"
$ws1.Range("N1").Copy()
$ws1.Range("L2").PasteSpecial(-4163)   # xlPasteValues
$ws1.Range("N1").Delete(-4159)         # xlShiftToLeft

# -----------------------------------------------------------------------
# Sheet "EditOverageDetails": update AdminUserName (F2) - the claims email
# used for admin notifications.
# -----------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("EditOverageDetails")
$ws4.Range("F2").Value = "info@freightclub.com"

# -----------------------------------------------------------------------
# Move the active tab / selection from "EditSecondaryInvoice" to
# "EditOverageDetails", with E2 selected there.
# -----------------------------------------------------------------------
$ws4.Activate()
$ws4.Range("E2").Select() | Out-Null
